# Add the new "SSR Relé Modul 2 kanál..." relay component as row 19 of the
# components list, extend the existing shared price-formula range down to
# E19, flag the single-channel relay (row 16) in red as a "duplicitně"-style
# warning, and tidy up the sheet view (selection/scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New component row (row 19): SSR relay, 2-channel -----------------
# Write the hyperlink target text first so it lands in the shared-string
# table ahead of the Czech label (matches the order new strings were
# appended upstream).
$ws.Range("F19").Value = "https://dratek.cz/arduino/1347-ssr-rele-modul-2-kanaly-5vdc-250vac-omron-g3mb-202p-solid-state-pro-arduino.html"
$ws.Range("B19").Value = "SSR Relé Modul 2 kanál 5VDC Low Level"
$ws.Range("C19").Value = 76
$ws.Range("D19").Value = 1
$ws.Range("E19").Formula = "=C19*D19"

# Turn F19 into a real hyperlink (adds the relationship / hyperlinks entry).
$ws.Hyperlinks.Add($ws.Range("F19"), "https://dratek.cz/arduino/1347-ssr-rele-modul-2-kanaly-5vdc-250vac-omron-g3mb-202p-solid-state-pro-arduino.html") | Out-Null

# Match formatting of the other rows in the same block: E column uses the
# bold "price" currency style, F column uses the hyperlink style. Grab the
# formats (only) from neighboring rows so we reuse the existing style
# entries instead of minting new ones.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The shared price formula that used to cover E8:E18 now needs to include
# the new row.
$ws.Range("E8:E19").FormulaR1C1 = "=RC[-2]*RC[-1]"

# --- Flag the 1-channel SSR relay (row 16) in red --------------------
$ws.Range("B16").Font.Color = 255

# --- Sheet view: scroll back to the top and select F9:F15 ------------
$ws.Range("F9:F15").Select()
